# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (sheet1) and "全部类型" (sheet4) worksheets to the refreshed
# values captured at the later gh-pages build (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item("展览")
$sheetAll  = $wb.Worksheets.Item("全部类型")

# Row => new F value, for the "展览" sheet
$expoUpdates = @{
    2  = 14877
    3  = 18506
    5  = 113
    7  = 220
    14 = 104
    15 = 198
    16 = 53
    17 = 1411
    20 = 84
    21 = 227
    22 = 7675
    24 = 20
    26 = 1216
    28 = 5954
    29 = 99
    30 = 62
    34 = 5301
}

foreach ($row in $expoUpdates.Keys) {
    $sheetExpo.Range("F$row").Value = $expoUpdates[$row]
}

# Row => new F value, for the "全部类型" sheet
$allUpdates = @{
    2  = 14877
    3  = 18506
    5  = 113
    7  = 220
    14 = 104
    15 = 198
    16 = 53
    17 = 1411
    21 = 84
    22 = 227
    23 = 7675
    25 = 20
    27 = 1217
    31 = 5954
    32 = 99
    33 = 62
    37 = 5301
}

foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
